$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 62.74008966666666
$ws.Range("H2").Value = 188.220269
$ws.Range("I2").Value = 0.6543216524118473
$ws.Range("J2").Value = 0.6543216524118471
$ws.Range("M2").Value = 3.135398666666667
$ws.Range("N2").Value = 9.406196000000001
$ws.Range("O2").Value = 0.1723049126704688
$ws.Range("P2").Value = 0.1723049126704688
$ws.Range("Q2").Value = 196.7151934874138
$ws.Range("R2").Value = 1770.436741386724
$ws.Range("S2").Value = 0.1127428351772202
$ws.Range("T2").Value = 0.1127428351772202
$ws.Range("G3").Value = 62.74008966666666
$ws.Range("H3").Value = 188.220269
$ws.Range("I3").Value = 0.6543216524118473
$ws.Range("J3").Value = 0.6543216524118471
$ws.Range("O3").Value = 0.1733096678828815
$ws.Range("P3").Value = 0.1733096678828815
$ws.Range("Q3").Value = 197.8622914601526
$ws.Range("R3").Value = 1780.760623141374
$ws.Range("S3").Value = 0.1134002682680754
$ws.Range("T3").Value = 0.1134002682680754
$ws.Range("G4").Value = 62.74008966666666
$ws.Range("H4").Value = 188.220269
$ws.Range("I4").Value = 0.6543216524118473
$ws.Range("J4").Value = 0.6543216524118471
$ws.Range("M4").Value = 0.4900660000000001
$ws.Range("N4").Value = 1.470198
$ws.Range("O4").Value = 0.02693143306797965
$ws.Range("P4").Value = 0.02693143306797965
$ws.Range("Q4").Value = 30.74678478258467
$ws.Range("R4").Value = 276.721063043262
$ws.Range("S4").Value = 0.01762181978685951
$ws.Range("T4").Value = 0.01762181978685951
$ws.Range("G5").Value = 62.74008966666666
$ws.Range("H5").Value = 188.220269
$ws.Range("I5").Value = 0.6543216524118473
$ws.Range("J5").Value = 0.6543216524118471
$ws.Range("M5").Value = 11.417657
$ws.Range("N5").Value = 34.252971
$ws.Range("O5").Value = 0.62745398637867
$ws.Range("P5").Value = 0.6274539863786701
$ws.Range("Q5").Value = 716.3448239632443
$ws.Range("R5").Value = 6447.103415669199
$ws.Range("S5").Value = 0.410556729179692
$ws.Range("T5").Value = 0.410556729179692
$ws.Range("I6").Value = 0.1782000513806195
$ws.Range("J6").Value = 0.1782000513806195
$ws.Range("M6").Value = 3.135398666666667
$ws.Range("N6").Value = 9.406196000000001
$ws.Range("O6").Value = 0.1723049126704688
$ws.Range("P6").Value = 0.1723049126704688
$ws.Range("Q6").Value = 53.57404490221779
$ws.Range("R6").Value = 482.1664041199601
$ws.Range("S6").Value = 0.03070474429101071
$ws.Range("T6").Value = 0.0307047442910107
$ws.Range("I7").Value = 0.1782000513806195
$ws.Range("J7").Value = 0.1782000513806195
$ws.Range("O7").Value = 0.1733096678828815
$ws.Range("P7").Value = 0.1733096678828815
$ws.Range("S7").Value = 0.03088379172148758
$ws.Range("T7").Value = 0.03088379172148758
$ws.Range("I8").Value = 0.1782000513806195
$ws.Range("J8").Value = 0.1782000513806195
$ws.Range("M8").Value = 0.4900660000000001
$ws.Range("N8").Value = 1.470198
$ws.Range("O8").Value = 0.02693143306797965
$ws.Range("P8").Value = 0.02693143306797965
$ws.Range("Q8").Value = 8.373677697886668
$ws.Range("R8").Value = 75.36309928098001
$ws.Range("S8").Value = 0.00479918275646769
$ws.Range("T8").Value = 0.004799182756467689
$ws.Range("I9").Value = 0.1782000513806195
$ws.Range("J9").Value = 0.1782000513806195
$ws.Range("M9").Value = 11.417657
$ws.Range("N9").Value = 34.252971
$ws.Range("O9").Value = 0.62745398637867
$ws.Range("P9").Value = 0.6274539863786701
$ws.Range("Q9").Value = 195.0916402750233
$ws.Range("R9").Value = 1755.82476247521
$ws.Range("S9").Value = 0.1118123326116535
$ws.Range("T9").Value = 0.1118123326116535
$ws.Range("G10").Value = 2.950144666666667
$ws.Range("H10").Value = 8.850434
$ws.Range("I10").Value = 0.03076730593473967
$ws.Range("J10").Value = 0.03076730593473966
$ws.Range("M10").Value = 3.135398666666667
$ws.Range("N10").Value = 9.406196000000001
$ws.Range("O10").Value = 0.1723049126704688
$ws.Range("P10").Value = 0.1723049126704688
$ws.Range("Q10").Value = 9.249879654340445
$ws.Range("R10").Value = 83.24891688906401
$ws.Range("S10").Value = 0.005301357962190915
$ws.Range("T10").Value = 0.005301357962190914
$ws.Range("G11").Value = 2.950144666666667
$ws.Range("H11").Value = 8.850434
$ws.Range("I11").Value = 0.03076730593473967
$ws.Range("J11").Value = 0.03076730593473966
$ws.Range("O11").Value = 0.1733096678828815
$ws.Range("P11").Value = 0.1733096678828815
$ws.Range("Q11").Value = 9.303818132662666
$ws.Range("R11").Value = 83.734363193964
$ws.Range("S11").Value = 0.005332271573200739
$ws.Range("T11").Value = 0.005332271573200738
$ws.Range("G12").Value = 2.950144666666667
$ws.Range("H12").Value = 8.850434
$ws.Range("I12").Value = 0.03076730593473967
$ws.Range("J12").Value = 0.03076730593473966
$ws.Range("M12").Value = 0.4900660000000001
$ws.Range("N12").Value = 1.470198
$ws.Range("O12").Value = 0.02693143306797965
$ws.Range("P12").Value = 0.02693143306797965
$ws.Range("Q12").Value = 1.445765596214667
$ws.Range("R12").Value = 13.011890365932
$ws.Range("S12").Value = 0.0008286076404634944
$ws.Range("T12").Value = 0.0008286076404634942
$ws.Range("G13").Value = 2.950144666666667
$ws.Range("H13").Value = 8.850434
$ws.Range("I13").Value = 0.03076730593473967
$ws.Range("J13").Value = 0.03076730593473966
$ws.Range("M13").Value = 11.417657
$ws.Range("N13").Value = 34.252971
$ws.Range("O13").Value = 0.62745398637867
$ws.Range("P13").Value = 0.6274539863786701
$ws.Range("Q13").Value = 33.68373990437934
$ws.Range("R13").Value = 303.153659139414
$ws.Range("S13").Value = 0.01930506875888452
$ws.Range("T13").Value = 0.01930506875888451
$ws.Range("G14").Value = 13.108629
$ws.Range("H14").Value = 39.325887
$ws.Range("I14").Value = 0.1367109902727936
$ws.Range("J14").Value = 0.1367109902727935
$ws.Range("M14").Value = 3.135398666666667
$ws.Range("N14").Value = 9.406196000000001
$ws.Range("O14").Value = 0.1723049126704688
$ws.Range("P14").Value = 0.1723049126704688
$ws.Range("Q14").Value = 41.10077788842801
$ws.Range("R14").Value = 369.9070009958521
$ws.Range("S14").Value = 0.02355597524004701
$ws.Range("T14").Value = 0.023555975240047
$ws.Range("G15").Value = 13.108629
$ws.Range("H15").Value = 39.325887
$ws.Range("I15").Value = 0.1367109902727936
$ws.Range("J15").Value = 0.1367109902727935
$ws.Range("O15").Value = 0.1733096678828815
$ws.Range("P15").Value = 0.1733096678828815
$ws.Range("Q15").Value = 41.340447321978
$ws.Range("R15").Value = 372.064025897802
$ws.Range("S15").Value = 0.02369333632011769
$ws.Range("T15").Value = 0.02369333632011769
$ws.Range("G16").Value = 13.108629
$ws.Range("H16").Value = 39.325887
$ws.Range("I16").Value = 0.1367109902727936
$ws.Range("J16").Value = 0.1367109902727935
$ws.Range("M16").Value = 0.4900660000000001
$ws.Range("N16").Value = 1.470198
$ws.Range("O16").Value = 0.02693143306797965
$ws.Range("P16").Value = 0.02693143306797965
$ws.Range("Q16").Value = 6.424093379514001
$ws.Range("R16").Value = 57.816840415626
$ws.Range("S16").Value = 0.003681822884188957
$ws.Range("T16").Value = 0.003681822884188956
$ws.Range("G17").Value = 13.108629
$ws.Range("H17").Value = 39.325887
$ws.Range("I17").Value = 0.1367109902727936
$ws.Range("J17").Value = 0.1367109902727935
$ws.Range("M17").Value = 11.417657
$ws.Range("N17").Value = 34.252971
$ws.Range("O17").Value = 0.62745398637867
$ws.Range("P17").Value = 0.6274539863786701
$ws.Range("Q17").Value = 149.669829662253
$ws.Range("R17").Value = 1347.028466960277
$ws.Range("S17").Value = 0.0857798558284399
$ws.Range("T17").Value = 0.0857798558284399
